# Daily attendance processing - 2025-11-21 07:44:35
# Swap the first and last comma-separated entries in the "Recorded By"
# column (column G) for every data row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ', '
        if ($parts.Count -ge 2) {
            $last = $parts.Count - 1
            $first = $parts[0]
            $parts[0] = $parts[$last]
            $parts[$last] = $first
            $cell.Value = [string]::Join(', ', $parts)
        }
    }
}
